# Apply the weekly price-sheet refresh: shift each row's data up one
# position (row 2 -> row 11, row 3 -> row 2, ..., row 11 -> row 10 -> ... )
# i.e. row r now shows what used to be at the mapped source row below.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44482
$ws.Range("M2").Value = 160
$ws.Range("N2").Value = 25000
$ws.Range("O2").Value = 26000
$ws.Range("P2").Value = 25500
$ws.Range("Q2").Value = '$/caja 12 kilos'
$ws.Range("S2").Value = 2125
$ws.Range("T2").Value = 12

# Row 3
$ws.Range("D3").Value = 44524
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 23000
$ws.Range("O3").Value = 24000
$ws.Range("P3").Value = 23500
$ws.Range("S3").Value = 1958

# Row 4
$ws.Range("D4").Value = 44160
$ws.Range("L4").Value = 'Segunda'
$ws.Range("N4").Value = 19000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 19500
$ws.Range("Q4").Value = '$/caja 13 kilos'
$ws.Range("S4").Value = 1500
$ws.Range("T4").Value = 13

# Row 5
$ws.Range("D5").Value = 44545
$ws.Range("L5").Value = 'Primera'
$ws.Range("N5").Value = 23000
$ws.Range("O5").Value = 24000
$ws.Range("P5").Value = 23500
$ws.Range("Q5").Value = '$/bandeja 12 kilos'
$ws.Range("S5").Value = 1958

# Row 6
$ws.Range("D6").Value = 44489
$ws.Range("L6").Value = 'Primera'
$ws.Range("N6").Value = 24000
$ws.Range("O6").Value = 25000
$ws.Range("P6").Value = 24500
$ws.Range("Q6").Value = '$/caja 12 kilos'
$ws.Range("S6").Value = 2042
$ws.Range("T6").Value = 12

# Row 7
$ws.Range("D7").Value = 44441
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 29000
$ws.Range("O7").Value = 30000
$ws.Range("P7").Value = 29500
$ws.Range("Q7").Value = '$/caja 12 kilos'
$ws.Range("S7").Value = 2458
$ws.Range("T7").Value = 12

# Row 8
$ws.Range("D8").Value = 44167
$ws.Range("L8").Value = 'Segunda'
$ws.Range("N8").Value = 18000
$ws.Range("O8").Value = 19000
$ws.Range("P8").Value = 18500
$ws.Range("Q8").Value = '$/caja 13 kilos'
$ws.Range("S8").Value = 1423
$ws.Range("T8").Value = 13

# Row 9
$ws.Range("D9").Value = 44475
$ws.Range("L9").Value = 'Especial'
$ws.Range("N9").Value = 32000
$ws.Range("O9").Value = 33000
$ws.Range("P9").Value = 32500
$ws.Range("S9").Value = 2708

# Row 10
$ws.Range("D10").Value = 44496
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 23000
$ws.Range("O10").Value = 24000
$ws.Range("P10").Value = 23500
$ws.Range("S10").Value = 1958

# Row 11
$ws.Range("D11").Value = 44468
$ws.Range("N11").Value = 29000
$ws.Range("O11").Value = 30000
$ws.Range("P11").Value = 29500
$ws.Range("Q11").Value = '$/bandeja 10 kilos'
$ws.Range("S11").Value = 2950
$ws.Range("T11").Value = 10
